# Daily attendance processing - 2026-01-30 22:40:06
# Normalizes the "Recorded By" (column G) lists so that the most recently
# recorded-by actor is listed first: the last comma-separated name/email in
# each cell is moved to the front of the list (rotate-right-by-1). Entries
# that include "admin@admin.com" are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    if (-not ($val.Contains(","))) { continue }
    if ($val.Contains("admin@admin.com")) { continue }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $lastIdx = $parts.Length - 1
    $lastItem = $parts[$lastIdx]
    $remaining = $parts[0..($lastIdx - 1)]

    $newParts = @($lastItem) + $remaining
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value2 = $newVal
}
